# Insert a new row at 239, shifting existing rows 239:329 down to 240:330,
# then populate the new row with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 239 (pushes rows 239-329 down to 240-330)
$ws.Rows.Item(239).Insert()

# Populate the newly inserted row 239 with the new record
$ws.Cells.Item(239, 1).Value = 4
$ws.Cells.Item(239, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(239, 3).Value = "Los Lagos"
$ws.Cells.Item(239, 4).Value = 44924
$ws.Cells.Item(239, 5).Value = 10
$ws.Cells.Item(239, 6).Value = 100112044
$ws.Cells.Item(239, 7).Value = "Perejil"
$ws.Cells.Item(239, 8).Value = "Sin especificar"
$ws.Cells.Item(239, 9).Value = "Primera"
$ws.Cells.Item(239, 10).Value = 60
$ws.Cells.Item(239, 11).Value = 6000
$ws.Cells.Item(239, 12).Value = 7000
$ws.Cells.Item(239, 13).Value = 6500
$ws.Cells.Item(239, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(239, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(239, 16).Value = 3250
$ws.Cells.Item(239, 17).Value = 2
$ws.Cells.Item(239, 18).Value = "Hortaliza"

# Ensure the date cell uses the same number format (style) as the rest of column D
$ws.Cells.Item(239, 4).NumberFormat = $ws.Cells.Item(240, 4).NumberFormat
